$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.388.63"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "3.556.16"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.68"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.53"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").Value = "3.555.81"
$ws.Range("E7").Value = "  +0.93%  "

$ws.Range("E9").Value = "  +2.47%  "

$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.86"
$ws.Range("E11").Value = "  -2.62%  "

$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("D13").Value = "4.160.01"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.15"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").Value = "3.552.54"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("D17").Value = "66.464.66"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.51"
$ws.Range("E19").Value = "  +5.75%  "

$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.86"
$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.39"
$ws.Range("E22").Value = "  +1.12%  "

$ws.Range("E23").Value = "  +1.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.61"
$ws.Range("E24").Value = "  +1.11%  "

$ws.Range("D25").Value = "3.697.12"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  +0.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.51"
$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("E29").Value = "  -0.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.97"
$ws.Range("E30").Value = "  -0.52%  "

$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("D33").Value = "3.550.95"
$ws.Range("E33").Value = "  +0.88%  "

$ws.Range("E34").Value = "  -1.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.155"
$ws.Range("E35").Value = "  -4.25%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.83"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("B37").Value = "USDe"
$ws.Range("C37").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("E38").Value = "  -1.40%  "

$ws.Range("E39").Value = "  -0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.96"
$ws.Range("E40").Value = "  +1.96%  "

$ws.Range("E41").Value = "  -0.74%  "

$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("E44").Value = "  +1.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.00"
$ws.Range("E45").Value = "  +1.62%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").Value = "  +4.60%  "

$ws.Range("E48").Value = "  -1.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.11"
$ws.Range("E49").Value = "  -3.35%  "

$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("E51").Value = "  +2.59%  "
